# Updated cryptos list (prices + 1h volume change) to match the latest
# coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '1.597.48'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").Value = "'211.38"
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").Value = "'0.0619"
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("D10").Value = "'19.51"
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("D11").Value = "'0.0840"
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").Value = '1.821.88'
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").Value = '1.599.49'
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D16").Value = "'65.06"
$ws.Range("E16").Value = '  +0.50%  '
$ws.Range("D17").Value = '26.645.72'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '0.0₃0748'
$ws.Range("E18").Value = '  +2.73%  '
$ws.Range("D19").Value = "'209.85"
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("E21").Value = '  +4.31%  '
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("D23").Value = "'2.32"
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("E24").Value = '  +0.90%  '
$ws.Range("D25").Value = "'143.23"
$ws.Range("E25").Value = '  -1.83%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").Value = "'7.12"
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("E28").Value = '  -1.02%  '
$ws.Range("D29").Value = "'15.34"
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("E30").Value = '  +1.91%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Value = "'3.25"
$ws.Range("E32").Value = '  +0.27%  '
$ws.Range("E33").Value = '  +0.52%  '
$ws.Range("D34").Value = '1.286.71'
$ws.Range("D35").Value = "'0.620"
$ws.Range("E35").Value = '  -5.75%  '
$ws.Range("D36").Value = "'2.44"
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("E39").Value = '  +17.60%  '
$ws.Range("D40").Value = "'0.826"
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("D41").Value = "'5.44"
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("D43").Value = "'0.783"
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("D44").Value = "'63.27"
$ws.Range("D45").Value = '1.734.45'
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").Value = "'91.18"
$ws.Range("E46").Value = '  +1.18%  '
$ws.Range("D47").Value = "'1.57"
$ws.Range("E47").Value = '  -2.76%  '
$ws.Range("D49").Value = "'0.0510"
$ws.Range("E49").Value = '  +0.85%  '
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").Value = "'7.35"
$ws.Range("E51").Value = '  -1.44%  '
